$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - only column A (the leading-space "use" column) changes
Set-TextCell "A2" " 67"

# Row 3 - full row update
Set-TextCell "A3" " 20"
$ws.Range("B3").Value = 70
$ws.Range("C3").Value = 1.11
$ws.Range("D3").Value = 204
$ws.Range("E3").Value = 1.06
$ws.Range("F3").Value = 125.6
$ws.Range("G3").Value = 0.79
$ws.Range("H3").Value = 0.73
$ws.Range("I3").Value = 0.34
$ws.Range("J3").Value = 0.09
$ws.Range("K3").Value = 0.03
$ws.Range("L3").Value = 51
$ws.Range("M3").Value = 48
$ws.Range("N3").Value = 24
$ws.Range("O3").Value = 6
$ws.Range("P3").Value = 2

# Row 4 - full row update
Set-TextCell "A4" " 13"
$ws.Range("B4").Value = 48
$ws.Range("C4").Value = 1.07
$ws.Range("D4").Value = 163
$ws.Range("E4").Value = 0.84
$ws.Range("F4").Value = 110.7
$ws.Range("G4").Value = 0.73
$ws.Range("H4").Value = 0.54
$ws.Range("I4").Value = 0.44
$ws.Range("J4").Value = 0.15
$ws.Range("K4").Value = 0.08
$ws.Range("L4").Value = 26
$ws.Range("M4").Value = 31
$ws.Range("N4").Value = 21
$ws.Range("O4").Value = 7
$ws.Range("P4").Value = 4
